# Refresh the cryptocurrency price/volume table with the latest scrape
# (GitHub Actions cron). All D/E (and occasionally B/C) cells in this
# sheet are stored as literal text, including values that look numeric
# (e.g. "594.56"), so writes force text entry with a leading apostrophe
# and then reset the cell style to "Normal" -- otherwise Excel leaves a
# "Text" number-format behind on any numeric-looking literal, which would
# add style metadata that should not be there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

Set-TextCell 'D2' "67.692.34"
Set-TextCell 'E2' "  +0.54%  "
Set-TextCell 'D3' "2.614.83"
Set-TextCell 'E3' "  -0.72%  "
Set-TextCell 'E4' "  +0.04%  "
Set-TextCell 'D5' "594.56"
Set-TextCell 'E5' "  -1.33%  "
Set-TextCell 'D6' "152.33"
Set-TextCell 'E6' "  -0.71%  "
Set-TextCell 'E7' "  +0.06%  "
Set-TextCell 'D8' "0.542"
Set-TextCell 'E8' "  -3.42%  "
Set-TextCell 'D9' "2.613.91"
Set-TextCell 'E9' "  -0.63%  "
Set-TextCell 'E10' "  +6.22%  "
Set-TextCell 'E11' "  -0.57%  "
Set-TextCell 'D12' "5.19"
Set-TextCell 'E12' "  -0.42%  "
Set-TextCell 'D13' "0.345"
Set-TextCell 'E13' "  -2.12%  "
Set-TextCell 'D14' "27.41"
Set-TextCell 'E14' "  -1.54%  "
Set-TextCell 'E15' "  +2.29%  "
Set-TextCell 'D16' "3.090.69"
Set-TextCell 'E16' "  -0.69%  "
Set-TextCell 'D17' "67.565.62"
Set-TextCell 'E17' "  +0.47%  "
Set-TextCell 'D18' "2.603.13"
Set-TextCell 'E18' "  -1.02%  "
Set-TextCell 'D19' "370.92"
Set-TextCell 'E19' "  +1.91%  "
Set-TextCell 'D20' "11.17"
Set-TextCell 'E20' "  -0.60%  "
Set-TextCell 'D21' "4.21"
Set-TextCell 'E21' "  -2.09%  "
Set-TextCell 'E22' "  -13.11%  "
Set-TextCell 'D23' "4.77"
Set-TextCell 'E23' "  -3.54%  "
Set-TextCell 'D24' "2.03"
Set-TextCell 'E24' "  -4.96%  "
Set-TextCell 'D25' "72.89"
Set-TextCell 'E25' "  +10.13%  "
Set-TextCell 'D26' "0.999"
Set-TextCell 'E26' "  -0.04%  "
Set-TextCell 'D27' "9.82"
Set-TextCell 'E27' "  -3.06%  "
Set-TextCell 'B28' "WrappedeETH"
Set-TextCell 'C28' "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextCell 'D28' "2.747.32"
Set-TextCell 'E28' "  -0.60%  "
Set-TextCell 'B29' "PEPE"
Set-TextCell 'C29' "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell 'D29' "0.0000103"
Set-TextCell 'E29' "  -0.56%  "
Set-TextCell 'D30' "587.40"
Set-TextCell 'E30' "  +1.51%  "
Set-TextCell 'E31' "  -0.42%  "
Set-TextCell 'E32' "  -2.22%  "
Set-TextCell 'D33' "7.77"
Set-TextCell 'E33' "  -1.53%  "
Set-TextCell 'D34' "1.83"
Set-TextCell 'E34' "  -0.92%  "
Set-TextCell 'E35' "  +0.07%  "
Set-TextCell 'E36' "  -2.33%  "
Set-TextCell 'E37' "  -1.88%  "
Set-TextCell 'D38' "158.75"
Set-TextCell 'E38' "  +0.83%  "
Set-TextCell 'D39' "19.12"
Set-TextCell 'E39' "  -1.57%  "
Set-TextCell 'D40' "1.88"
Set-TextCell 'E40' "  +2.85%  "
Set-TextCell 'D41' "0.366"
Set-TextCell 'E41' "  -1.08%  "
Set-TextCell 'D42' "5.26"
Set-TextCell 'E42' "  -0.38%  "
Set-TextCell 'D43' "2.64"
Set-TextCell 'E43' "  +0.43%  "
Set-TextCell 'D44' "17.08"
Set-TextCell 'E44' "  +4.47%  "
Set-TextCell 'E45' "  +0.04%  "
Set-TextCell 'D46' "40.39"
Set-TextCell 'E46' "  -1.97%  "
Set-TextCell 'D47' "0.0$([char]8326)0300"
Set-TextCell 'E47' "  +4.74%  "
Set-TextCell 'D48' "154.67"
Set-TextCell 'E48' "  -0.91%  "
Set-TextCell 'D49' "3.67"
Set-TextCell 'E49' "  -1.82%  "
Set-TextCell 'D50' "1.68"
Set-TextCell 'E50' "  -2.54%  "
Set-TextCell 'D51' "0.0776"
Set-TextCell 'E51' "  -1.97%  "
